# WAT new test cases
# Adds a new test case row (WAT23 / WAT-142) to the "Test Cases" sheet,
# matching the formatting of the row immediately above it, with the
# description cell wrapped and a taller row to fit it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Clone the formatting of the last existing row (23) onto the new row (24)
# so the new row picks up the same borders/fills as the rest of the table.
$ws.Range("A23:E23").Copy()
$ws.Range("A24:E24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new test case. Order matches the TCID / Description / JIRA ID /
# Runmode column layout used throughout the sheet.
$ws.Range("A24").Value = "WAT23"
$ws.Range("C24").Value = "Verify that system must display following publication details in cart if publication details 1, Publication count 1, Year, journal, published date, author metadata ..etc."
$ws.Range("B24").Value = "WAT-142"
$ws.Range("D24").Value = "Y"

# Wrap the long description and give the row enough height to show it.
$ws.Range("C24").WrapText = $true
$ws.Rows.Item(24).RowHeight = 30

# Move the sheet's active selection down past the new row, as in the source file.
$ws.Range("C29").Select() | Out-Null
